$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.3921007115659704
$ws.Range("C3").Value = 0.2244505373841501
$ws.Range("C4").Value = 0.5373564881718573
$ws.Range("C5").Value = 0.4391317600812851
$ws.Range("C6").Value = 0.5602820546308499
$ws.Range("C7").Value = 0.4281684744334564
$ws.Range("C8").Value = 0.5188885972004412
$ws.Range("C9").Value = 0.2406766632537068
$ws.Range("C10").Value = 0.2969811255454097
